$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-13 as recalculated in the regen.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 8
$ws.Range("G7").Value = 4
$ws.Range("G9").Value = 3
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 3
$ws.Range("G12").Value = 4
$ws.Range("G13").Value = 2
